# Append a freshly-scraped Lancers listing to the "ランサーズ" sheet and
# refresh the "取得日時" (fetched-at) timestamp on every existing row to
# match this scrape run: 2025-10-03 18:23:27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newStamp = "2025-10-03 18:23:27"

# Header is row 1; existing data currently runs through the last used row.
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Update column A (取得日時) for every already-present data row (2..lastRow).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newStamp
}

# Append the newly scraped row right after the existing data.
$ws.Cells.Item($newRow, 1).Value = $newStamp
$ws.Cells.Item($newRow, 2).Value = "【個人利用】Web情報収集の仕組み構築ご依頼"
$ws.Cells.Item($newRow, 3).Value = "システム開発"
$ws.Cells.Item($newRow, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item($newRow, 5).Value = "期限情報なし"

$url = "https://www.lancers.jp/work/detail/5406440"
$ws.Cells.Item($newRow, 6).Value = $url
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 6), $url)
$ws.Cells.Item($newRow, 6).Style = "Hyperlink"

$ws.Cells.Item($newRow, 7).Value = 10

"appended row $newRow"
